$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.366.48'
$ws.Range("E2").Value = '  +1.78%  '

$ws.Range("D3").Value = '1.882.27'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '246.63'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.58%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.28'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +4.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.356'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0745'
$ws.Range("D10").ClearFormats()

$ws.Range("E11").Value = '  +1.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '13.50'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.86%  '

$ws.Range("D13").Value = '2.158.60'
$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.771'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.92%  '

$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("D16").Value = '1.896.73'
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("D17").Value = '35.342.52'
$ws.Range("E17").Value = '  +1.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '73.42'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.79%  '

$ws.Range("E19").Value = '  +0.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '244.63'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.82'
$ws.Range("D21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.19'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +5.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.60'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +8.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("E25").Value = '  -3.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.57'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.64'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.30'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("E30").Value = '  +2.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.28'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.88'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.18'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.48%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("E35").Value = '  -12.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.853'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.38%  '

$ws.Range("E37").Value = '  -1.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0734'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +11.63%  '

$ws.Range("E39").Value = '  +0.66%  '

$ws.Range("E40").Value = '  +3.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.25'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.87%  '

$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("E43").Value = '  +1.92%  '

$ws.Range("D44").Value = '1.308.28'
$ws.Range("E44").Value = '  +1.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0813'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.73'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.09'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.32%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.33'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.28'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.65%  '

$ws.Range("D51").Value = '2.062.18'
$ws.Range("E51").Value = '  +0.24%  '
